# "Generate Report for handback"
#
# Two source files (a1069261-....md and its dependent d4902239-....md) have
# now been handed back (both zh-cn and de-de). This updates:
#   - Overview sheet: Status columns for those two rows -> "Handed back: in sync with en-US"
#   - zh-cn / de-de sheets: Status -> handed back, plus the "Latest Target File",
#     "Latest Handback File" and "Latest Handback DateTime" columns are now
#     populated for those two rows (they previously were blank / "never").

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: rows for a1069261-....md (row 3) and d4902239-....md (row 4)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack
$overview.Range("B4").Value = $handedBack
$overview.Range("C4").Value = $handedBack

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Row 3 - a1069261-....md
$zhcn.Range("B3").Value = $handedBack
$zhcn.Range("E3").Value = "a1069261-3c07-4184-9a12-11e9257960f2.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/40c56059b388dc319f551e5bc3caece3e0ee0bdb/e2e/a1069261-3c07-4184-9a12-11e9257960f2.md", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.md") | Out-Null
$zhcn.Range("F3").Value = "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3c6622ee8e332818ac3b7f750e234e2e46da4541/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.zh-cn.xlf", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.zh-cn.xlf") | Out-Null
$zhcn.Range("G3").Value = "2016-01-18 06:42:26"

# Row 4 - d4902239-....md (depends on a1069261)
$zhcn.Range("B4").Value = $handedBack
$zhcn.Range("E4").Value = "a1069261-3c07-4184-9a12-11e9257960f2.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/40c56059b388dc319f551e5bc3caece3e0ee0bdb/e2e/a1069261-3c07-4184-9a12-11e9257960f2.md", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.md") | Out-Null
$zhcn.Range("F4").Value = "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3c6622ee8e332818ac3b7f750e234e2e46da4541/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.zh-cn.xlf", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.zh-cn.xlf") | Out-Null
$zhcn.Range("G4").Value = "2016-01-18 06:42:26"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Row 3 - a1069261-....md
$dede.Range("B3").Value = $handedBack
$dede.Range("E3").Value = "a1069261-3c07-4184-9a12-11e9257960f2.md"
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/40c56059b388dc319f551e5bc3caece3e0ee0bdb/e2e/a1069261-3c07-4184-9a12-11e9257960f2.md", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.md") | Out-Null
$dede.Range("F3").Value = "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/df9a73f113d603d669fa1dd5750f8977ff56435c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.de-de.xlf", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.de-de.xlf") | Out-Null
$dede.Range("G3").Value = "2016-01-18 06:42:42"

# Row 4 - d4902239-....md (depends on a1069261)
$dede.Range("B4").Value = $handedBack
$dede.Range("E4").Value = "a1069261-3c07-4184-9a12-11e9257960f2.md"
$dede.Hyperlinks.Add($dede.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/40c56059b388dc319f551e5bc3caece3e0ee0bdb/e2e/a1069261-3c07-4184-9a12-11e9257960f2.md", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.md") | Out-Null
$dede.Range("F4").Value = "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/df9a73f113d603d669fa1dd5750f8977ff56435c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.de-de.xlf", "", "", "a1069261-3c07-4184-9a12-11e9257960f2.49ff2e76304917b09bbec05a8384ffce17f3f250.de-de.xlf") | Out-Null
$dede.Range("G4").Value = "2016-01-18 06:42:42"
